$d = $word.ActiveDocument

# The last paragraph in the document contains a single run with text "aa".
# Remove that run's text, leaving the paragraph empty (matching the diff).
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.End = $r.End - 1
$r.Text = ""
